# Apply the edits to the EC_Prism_Template worksheet (active sheet):
# - Rename "EC Exists" -> "EC Exists with mu=0" (F2)
# - Rename "Prob for EC" -> "Prob for EC with mu=0" (G2)
# - Add two new columns: "EC Exists with mu" (H2), "Prob for EC with mu" (I2)
# - Extend header merge A1:G1 -> A1:I1
# - Adjust column widths
# - Move the selection to J12

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EC_Prism_Template")

# Update existing header text
$ws.Range("F2").Value = "EC Exists with mu=0"
$ws.Range("G2").Value = "Prob for EC with mu=0"

# Copy formatting of existing header cells (F2/G2) onto the new header cells
$ws.Range("G2").Copy()
$ws.Range("H2:I2").PasteSpecial(-4122)
$ws.Range("G3").Copy()
$ws.Range("H3:I3").PasteSpecial(-4122)
$ws.Range("G1").Copy()
$ws.Range("H1:I1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Set new header values
$ws.Range("H2").Value = "EC Exists with mu"
$ws.Range("I2").Value = "Prob for EC with mu"

# Re-merge header row across the new range
$ws.Range("A1:G1").UnMerge()
$ws.Range("A1:I1").Merge()

# Merge the two new header cells vertically with row 3
$ws.Range("H2:H3").Merge()
$ws.Range("I2:I3").Merge()

# Column widths to match the new layout
$ws.Columns.Item("F").ColumnWidth = 10.91
$ws.Columns.Item("G").ColumnWidth = 13.67
$ws.Columns.Item("H").ColumnWidth = 8.18

# Update selection to match target workbook state
$ws.Range("J12").Select()
